$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking strings such as
# "64.560.33", "0.999" or "0.0000258". Excel would normally parse
# these as numbers (or reject "64.560.33" outright) when assigned
# through .Value, which would not match the original inline-string
# cell contents. Setting NumberFormat to Text ("@") first forces
# Excel to keep the assigned value as literal text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.560.33"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.132.94"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.98"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.67"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.133.15"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.04"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.493"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +7.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.637.24"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.711.44"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.131.70"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.04"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "498.33"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.06"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.63"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.38"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.80"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.87"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.31"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.12"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.40"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0893"
$ws.Range("E37").Value = "  +6.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "460.36"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0414"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.57"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.018.20"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.281"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.05"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0567"
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.90"
$ws.Range("E51").Value = "  +0.14%  "
